$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G1:G48").Copy($ws.Range("H1:H48"))
